# Apply the Aug 17 2023 17:50:32 UTC cryptos-list refresh (GitHub Actions price/volume pull).
# All affected cells in this sheet are stored as text (t="inlineStr"/shared string), including
# Price values that look numeric (e.g. "1.002", "0.4422") and multi-dot strings (e.g. "28.019.43").
# Plain `.Value = "1.002"` would get auto-coerced to the number 1.002 (losing trailing zeros / dot
# grouping), so each write forces text via NumberFormat "@" first, then clears the now-unneeded
# number format back off the cell (ClearFormats) to avoid leaving a stray style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '28.019.43'
Set-TextValue 'E2' '  -3.98%  '
Set-TextValue 'D3' '1.742.13'
Set-TextValue 'E3' '  -4.61%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '226.38'
Set-TextValue 'E5' '  -3.48%  '
Set-TextValue 'D6' '0.5780'
Set-TextValue 'E6' '  -3.77%  '
Set-TextValue 'D7' '1.003'
Set-TextValue 'E7' '  -0.03%  '
Set-TextValue 'D8' '0.2710'
Set-TextValue 'E8' '  -1.49%  '
Set-TextValue 'D9' '23.12'
Set-TextValue 'E9' '  -1.31%  '
Set-TextValue 'D10' '0.06582'
Set-TextValue 'E10' '  -5.19%  '
Set-TextValue 'D11' '0.07535'
Set-TextValue 'E11' '  -0.65%  '
Set-TextValue 'D12' '1.737.83'
Set-TextValue 'E12' '  -5.02%  '
Set-TextValue 'D13' '4.711'
Set-TextValue 'E13' '  -0.44%  '
Set-TextValue 'D14' '0.6043'
Set-TextValue 'E14' '  -3.25%  '
Set-TextValue 'D15' '1.980.04'
Set-TextValue 'E15' '  -4.58%  '
Set-TextValue 'D16' '74.26'
Set-TextValue 'E16' '  -4.06%  '
Set-TextValue 'D17' '0.000008672'
Set-TextValue 'E17' '  -11.30%  '
Set-TextValue 'D18' '28.021.16'
Set-TextValue 'E18' '  -2.99%  '
Set-TextValue 'D19' '5.320'
Set-TextValue 'E19' '  -4.49%  '
Set-TextValue 'E20' '  -0.11%  '
Set-TextValue 'D21' '205.06'
Set-TextValue 'E21' '  -5.31%  '
Set-TextValue 'D22' '11.26'
Set-TextValue 'E22' '  -2.37%  '
Set-TextValue 'D23' '6.623'
Set-TextValue 'E23' '  -3.88%  '
Set-TextValue 'E24' '  -0.05%  '
Set-TextValue 'D25' '149.83'
Set-TextValue 'E25' '  -3.97%  '
Set-TextValue 'D26' '8.053'
Set-TextValue 'E26' '  +1.60%  '
Set-TextValue 'D27' '0.1231'
Set-TextValue 'E27' '  -4.52%  '
Set-TextValue 'D28' '16.11'
Set-TextValue 'E28' '  -2.19%  '
Set-TextValue 'D29' '0.06179'
Set-TextValue 'E29' '  -5.76%  '
Set-TextValue 'D30' '1.382'
Set-TextValue 'E30' '  -1.87%  '
Set-TextValue 'D31' '1.391'
Set-TextValue 'E31' '  -3.32%  '
Set-TextValue 'D32' '3.739'
Set-TextValue 'E32' '  -2.18%  '
Set-TextValue 'D33' '3.710'
Set-TextValue 'E33' '  -1.64%  '
Set-TextValue 'D34' '1.673'
Set-TextValue 'E34' '  -2.91%  '
Set-TextValue 'D35' '1.034'
Set-TextValue 'E35' '  -5.34%  '
Set-TextValue 'D36' '0.6369'
Set-TextValue 'E36' '  -1.41%  '
Set-TextValue 'D37' '2.422'
Set-TextValue 'E37' '  -4.43%  '
Set-TextValue 'D38' '2.727'
Set-TextValue 'E38' '  -0.68%  '
Set-TextValue 'D39' '0.01669'
Set-TextValue 'E39' '  -5.06%  '
Set-TextValue 'D40' '1.128.08'
Set-TextValue 'E40' '  -1.32%  '
Set-TextValue 'D41' '6.197'
Set-TextValue 'E41' '  -4.00%  '
Set-TextValue 'D42' '0.8722'
Set-TextValue 'E42' '  -1.68%  '
Set-TextValue 'D43' '1.003'
Set-TextValue 'E43' '  +0.14%  '
Set-TextValue 'D44' '99.57'
Set-TextValue 'E44' '  -0.68%  '
Set-TextValue 'D45' '1.893.35'
Set-TextValue 'E45' '  -4.88%  '
Set-TextValue 'D46' '59.33'
Set-TextValue 'E46' '  -3.83%  '
Set-TextValue 'B47' 'BabyDogeCoin'
Set-TextValue 'C47' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D47' '0.00000000109'
Set-TextValue 'E47' '  -3.16%  '
Set-TextValue 'B48' 'RenderToken'
Set-TextValue 'C48' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D48' '1.576'
Set-TextValue 'E48' '  -2.15%  '
Set-TextValue 'D49' '8.269'
Set-TextValue 'E49' '  -2.54%  '
Set-TextValue 'D50' '0.05378'
Set-TextValue 'E50' '  -2.22%  '
Set-TextValue 'B51' 'Mantle'
Set-TextValue 'C51' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D51' '0.4422'
Set-TextValue 'E51' '  -2.51%  '
